# Lab 11 Traceroute and Path Finding - edits
$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Remove the stray _GoBack bookmark that sits after "Current IP: 10.3.40.113"
# ---------------------------------------------------------------------
try {
    $oldGoBack = $d.Bookmarks.Item("_GoBack")
    $oldGoBack.Delete()
} catch {
}

# ---------------------------------------------------------------------
# 2) Re-flow the six traceroute-hop paragraphs so each one's text is a
#    single contiguous run (collapsing the old keystroke-by-keystroke
#    run splits), except where the target keeps an explicit two-run
#    split (paragraph 2 "... - " / "Private Router." and the trailing
#    "." on paragraphs 5 and 6).
# ---------------------------------------------------------------------

$d.Content.Find.Execute("111.68.101.1 " + [char]0x2013 + " Islamabad Pakistan.", $true, $false, $false, $false, $false, $true, 1, $false, "111.68.101.1 " + [char]0x2013 + " Islamabad Pakistan.", 2) | Out-Null

$d.Content.Find.Execute("172.31.254.25 " + [char]0x2013 + " Private Router.", $true, $false, $false, $false, $false, $true, 1, $false, "172.31.254.25 " + [char]0x2013 + " Private Router.", 2) | Out-Null

$d.Content.Find.Execute("202.179.249.46 " + [char]0x2013 + " China Beijing.", $true, $false, $false, $false, $false, $true, 1, $false, "202.179.249.46 " + [char]0x2013 + " China Beijing.", 2) | Out-Null

$d.Content.Find.Execute("202.179.249.45 " + [char]0x2013 + " China Beijing.", $true, $false, $false, $false, $false, $true, 1, $false, "202.179.249.45 " + [char]0x2013 + " China Beijing.", 2) | Out-Null

$d.Content.Find.Execute("202.179.249.42 " + [char]0x2013 + " China Beijing.", $true, $false, $false, $false, $false, $true, 1, $false, "202.179.249.42 " + [char]0x2013 + " China Beijing.", 2) | Out-Null

$d.Content.Find.Execute("202.179.249.62 " + [char]0x2013 + " China Beijing.", $true, $false, $false, $false, $false, $true, 1, $false, "202.179.249.62 " + [char]0x2013 + " China Beijing.", 2) | Out-Null

# Now force a run boundary in paragraph 2 right before "Private Router."
# by round-tripping a formatting property (bold on/off) without leaving
# any visible trace - this produces a clean extra <w:r> split matching
# the target XML.
$rng = $d.Content
$found = $rng.Find.Execute("Private Router.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $rng.Font.Bold = 1
    $rng.Font.Bold = 0
}

# Paragraph 5: split off the trailing "." after "... China Beijing"
$rng = $d.Content
$found = $rng.Find.Execute("202.179.249.42 " + [char]0x2013 + " China Beijing.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $period = $d.Range($rng.End - 1, $rng.End)
    $period.Font.Bold = 1
    $period.Font.Bold = 0
}

# Paragraph 6: split off the trailing "." after "... China Beijing"
$rng = $d.Content
$found = $rng.Find.Execute("202.179.249.62 " + [char]0x2013 + " China Beijing.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $period = $d.Range($rng.End - 1, $rng.End)
    $period.Font.Bold = 1
    $period.Font.Bold = 0
}

# ---------------------------------------------------------------------
# 3) Rewrite the "Yes, it does..." sentence to "No it does not follow
#    same path, although it does in my two tries..." and re-home the
#    _GoBack bookmark to the split point between the two new runs.
# ---------------------------------------------------------------------
$d.Content.Find.Execute("Yes, it does in my two tries", $true, $false, $false, $false, $false, $true, 1, $false, "No it does not follow same path, although it does in my two tries", 2) | Out-Null

$rng = $d.Content
$found = $rng.Find.Execute("No it does not follow same path, although", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $splitPoint = $rng.End
    $whole = $d.Range($rng.Start, $rng.End)
    $whole.Font.Bold = 1
    $whole.Font.Bold = 0
    $d.Bookmarks.Add("_GoBack", $d.Range($splitPoint, $splitPoint)) | Out-Null
}
